# This workbook tracks weekly wholesale produce price reports.
# A new week of data (fecha = 44714) is being added at the top of the
# "Vega Modelo de Temuco - Cebolla" block, pushing the existing rows
# (962-1023) down by three rows (to 965-1026), and the dimension grows
# from A1:R1023 to A1:R1026.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 962; this shifts rows 962:1023 down to 965:1026
# and copies formatting (e.g. the date number format on column D) from the
# row above, matching how the existing rows were already formatted.
$ws.Rows("962:964").Insert()

# Row 962 - new price report entry
$ws.Range("A962").Value = 10
$ws.Range("B962").Value = "Vega Modelo de Temuco"
$ws.Range("C962").Value = "La Araucanía"
$ws.Range("D962").Value = 44714
$ws.Range("E962").Value = 9
$ws.Range("F962").Value = 100112004
$ws.Range("G962").Value = "Cebolla"
$ws.Range("H962").Value = "Morada(o)"
$ws.Range("I962").Value = "1a (guarda)"
$ws.Range("J962").Value = 250
$ws.Range("K962").Value = 13000
$ws.Range("L962").Value = 13000
$ws.Range("M962").Value = 13000
$ws.Range("N962").Value = "`$/malla 18 kilos"
$ws.Range("O962").Value = "Región de O'Higgins"
$ws.Range("P962").Value = 722
$ws.Range("Q962").Value = 18
$ws.Range("R962").Value = "Hortaliza"

# Row 963 - new price report entry
$ws.Range("A963").Value = 10
$ws.Range("B963").Value = "Vega Modelo de Temuco"
$ws.Range("C963").Value = "La Araucanía"
$ws.Range("D963").Value = 44714
$ws.Range("E963").Value = 9
$ws.Range("F963").Value = 100112004
$ws.Range("G963").Value = "Cebolla"
$ws.Range("H963").Value = "Sin especificar"
$ws.Range("I963").Value = "1a (guarda)"
$ws.Range("J963").Value = 650
$ws.Range("K963").Value = 6000
$ws.Range("L963").Value = 6000
$ws.Range("M963").Value = 6000
$ws.Range("N963").Value = "`$/malla 18 kilos"
$ws.Range("O963").Value = "Región del Maule"
$ws.Range("P963").Value = 333
$ws.Range("Q963").Value = 18
$ws.Range("R963").Value = "Hortaliza"

# Row 964 - new price report entry
$ws.Range("A964").Value = 10
$ws.Range("B964").Value = "Vega Modelo de Temuco"
$ws.Range("C964").Value = "La Araucanía"
$ws.Range("D964").Value = 44714
$ws.Range("E964").Value = 9
$ws.Range("F964").Value = 100112004
$ws.Range("G964").Value = "Cebolla"
$ws.Range("H964").Value = "Sin especificar"
$ws.Range("I964").Value = "1a (guarda)"
$ws.Range("J964").Value = 185
$ws.Range("K964").Value = 10000
$ws.Range("L964").Value = 10000
$ws.Range("M964").Value = 10000
$ws.Range("N964").Value = "`$/malla 25 kilos"
$ws.Range("O964").Value = "Región del Maule"
$ws.Range("P964").Value = 400
$ws.Range("Q964").Value = 25
$ws.Range("R964").Value = "Hortaliza"
